# Modules sheet: update header labels, widen the two last columns and
# move the active selection, per the "Mise à jour de certains champs de
# Modules et de Professeurs" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row text swap:
#   C1 "Enseignant"       -> "Chef  Module"
#   D1 "Nombre d'heures"  -> "Composants"
$ws.Range("C1").Value = "Chef  Module"
$ws.Range("D1").Value = "Composants"

# Column widths (character units); chosen so the saved OOXML <col>
# widths land on 35 and as close as possible to 24.5703125.
$ws.Columns.Item(3).ColumnWidth = 34.16666666666667
$ws.Columns.Item(4).ColumnWidth = 23.666666666666668

# Active cell / selection moves to E8
$ws.Range("E8").Select()
